$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("A18").Value = "No Emulator Data"
$ws.Range("C18").Value = "17"
$ws.Range("D18").Value = "3.0"
$ws.Range("E18").Value = "10.50"
$ws.Range("G18").Value = "PayNow"
$ws.Range("H18").Value = "en_US"
$ws.Range("W18").Value = "udf data 4"
$ws.Range("X18").Value = "udf data 5"
$ws.Range("AB18").Value = "udf data 9"
$ws.Range("AC18").Value = "udf data 10"

# Row 19
$ws.Range("A19").Value = "No Emulator Data"
$ws.Range("C19").Value = "18"
$ws.Range("D19").Value = "3.0"
$ws.Range("E19").Value = "10.50"
$ws.Range("G19").Value = "AutoPay"
$ws.Range("H19").Value = "en_US"
$ws.Range("W19").Value = "udf data 4"
$ws.Range("X19").Value = "udf data 5"
$ws.Range("AB19").Value = "udf data 9"
$ws.Range("AC19").Value = "udf data 10"

$ws.Range("D19").Select()
